$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(1).Insert()

$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "image"
$ws.Range("C1").Value = "description"

$null = $ws.Range("D1").Select()
